$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously-last row (row 6) had its Date cell (B6) tagged with the
# plain "date only" number format; the refresh re-tags it with the
# "date + time" format (style used by B2:B5) now that it is no longer the
# newest row.
$ws.Range("B6").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's fuel-price observation as row 7.
$ws.Range("A7").Value = 808.9640000000001
$ws.Range("B7").Value = 45735
$ws.Range("B7").NumberFormat = "YYYY-MM-DD"
$ws.Range("C7").Value = 806.651
